# Add carjacking data for 2022-03-29 (rolling "through" date moves from
# March 20 to March 21) - update sheet title, column header, and the
# affected neighborhood counts for the current month column (B) plus a
# handful of other cells that changed in the source data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Rename the sheet tab to reflect the new "through" date.
$ws.Name = "Through 2022-03-21"

# 2) Update the column header text (B1) to match.
$ws.Range("B1").Value = "March 2022 (through March 21)"

# 3) Row 3 - Austin
$ws.Range("B3").Value = 8
$ws.Range("E3").Value = 6

# 4) Row 5 - Garfield Park
$ws.Range("Q5").Value = 8

# 5) Row 13 - Woodlawn (new data point)
$ws.Range("W13").Value = 1

# 6) Row 15 - Humboldt Park
$ws.Range("B15").Value = 5

# 7) Row 25 - Washington Park (new data point)
$ws.Range("E25").Value = 1

# 8) Row 26 - Grand Crossing
$ws.Range("E26").Value = 4

# 9) Row 40 - Morgan Park (new data point)
$ws.Range("T40").Value = 1

# 10) Row 62 - Gage Park
$ws.Range("W62").Value = 2

# 11) Row 64 - Garfield Ridge (new data point)
$ws.Range("Q64").Value = 1
